$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2199
$ws.Range("J40").Value = 1766.6
$ws.Range("L40").Value = 1766.6
$ws.Range("N40").Value = -2116.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2729.5
$ws.Range("I64").Value = 2520.4
$ws.Range("K64").Value = 2520.4
$ws.Range("M64").Value = -2272.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2729.5
$ws.Range("I67").Value = 2520.4
$ws.Range("K67").Value = 2520.4
$ws.Range("M67").Value = -1662.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4430.5415
$ws.Range("I74").Value = 5350.3
$ws.Range("J74").Value = 3773.5715
$ws.Range("K74").Value = 5350.3
$ws.Range("L74").Value = 3773.5715
$ws.Range("M74").Value = -4414.3
$ws.Range("N74").Value = -5645.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 65281.438
$ws.Range("I76").Value = 74035.92999999999
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 74035.92999999999
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -73720.92999999999
$ws.Range("N76").Value = -4630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4430.5415
$ws.Range("I77").Value = 5350.3
$ws.Range("J77").Value = 3773.5715
$ws.Range("K77").Value = 26751.5
$ws.Range("L77").Value = 18867.8575
$ws.Range("M77").Value = -22071.5
$ws.Range("N77").Value = -28227.8575

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 65281.438
$ws.Range("I79").Value = 74035.92999999999
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 74035.92999999999
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -72943.92999999999
$ws.Range("N79").Value = -6184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1230.3928
$ws.Range("I137").Value = 762
$ws.Range("K137").Value = 2286
$ws.Range("M137").Value = 264

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2894.049
$ws.Range("I138").Value = 1898.75
$ws.Range("J138").Value = 3247.9333
$ws.Range("K138").Value = 5696.25
$ws.Range("L138").Value = 9743.7999
$ws.Range("M138").Value = -556.25
$ws.Range("N138").Value = -20023.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 58824772
$ws.Range("I2").Value = 250000480
$ws.Range("J2").Value = 1477.7693
$ws.Range("K2").Value = 250000480
$ws.Range("L2").Value = 1477.7693
$ws.Range("M2").Value = -250000367
$ws.Range("N2").Value = -1703.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 17544768
$ws.Range("I45").Value = 33334038
$ws.Range("J45").Value = 1134.8889
$ws.Range("K45").Value = 33334038
$ws.Range("L45").Value = 1134.8889
$ws.Range("M45").Value = -33333661
$ws.Range("N45").Value = -1888.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1103.6666
$ws.Range("I63").Value = 905.5
$ws.Range("K63").Value = 905.5
$ws.Range("M63").Value = -219.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1103.6666
$ws.Range("I66").Value = 905.5
$ws.Range("K66").Value = 4527.5
$ws.Range("M66").Value = -1095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 823.7406999999999
$ws.Range("I74").Value = 822.3019
$ws.Range("K74").Value = 822.3019
$ws.Range("M74").Value = 51.69809999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 823.7406999999999
$ws.Range("I77").Value = 822.3019
$ws.Range("K77").Value = 4111.5095
$ws.Range("M77").Value = 256.4904999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 58824772
$ws.Range("I116").Value = 250000480
$ws.Range("J116").Value = 1477.7693
$ws.Range("K116").Value = 250000480
$ws.Range("L116").Value = 1477.7693
$ws.Range("M116").Value = -249998186
$ws.Range("N116").Value = -6065.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8370.264999999999
$ws.Range("I132").Value = 9771.615
$ws.Range("J132").Value = 3815.875
$ws.Range("K132").Value = 29314.845
$ws.Range("L132").Value = 11447.625
$ws.Range("M132").Value = -26784.845
$ws.Range("N132").Value = -16507.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 58824772
$ws.Range("I3").Value = 250000480
$ws.Range("J3").Value = 1477.7693
$ws.Range("K3").Value = 250000480
$ws.Range("L3").Value = 1477.7693
$ws.Range("M3").Value = -250000366
$ws.Range("N3").Value = -1705.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3154.6365
$ws.Range("I105").Value = 2215
$ws.Range("J105").Value = 4282.2
$ws.Range("K105").Value = 2215
$ws.Range("L105").Value = 4282.2
$ws.Range("M105").Value = -468
$ws.Range("N105").Value = -7776.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1212.6666
$ws.Range("I58").Value = 1212.6666
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1212.6666
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = -1009.6666
$ws.Range("M58").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1212.6666
$ws.Range("I136").Value = 1212.6666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3637.9998
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -1087.9998
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 83337590
$ws.Range("J70").Value = 4502
$ws.Range("K70").Value = 83337590
$ws.Range("L70").Value = 4502
$ws.Range("M70").Value = -83337320
$ws.Range("N70").Value = -5042

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I73").Value = 83337590
$ws.Range("J73").Value = 4502
$ws.Range("K73").Value = 83337590
$ws.Range("L73").Value = 4502
$ws.Range("M73").Value = -83336654
$ws.Range("N73").Value = -6374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4343.5386
$ws.Range("I132").Value = 4606.25
$ws.Range("J132").Value = 3142.5715
$ws.Range("K132").Value = 13818.75
$ws.Range("L132").Value = 9427.7145
$ws.Range("M132").Value = -11288.75
$ws.Range("N132").Value = -14487.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33335726
$ws.Range("I7").Value = 2262.5
$ws.Range("J7").Value = 71431110
$ws.Range("K7").Value = 2262.5
$ws.Range("L7").Value = 71431110
$ws.Range("M7").Value = -2150.5
$ws.Range("N7").Value = -71431334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 33335726
$ws.Range("I126").Value = 2262.5
$ws.Range("J126").Value = 71431110
$ws.Range("K126").Value = 6787.5
$ws.Range("L126").Value = 214293330
$ws.Range("M126").Value = -4317.5
$ws.Range("N126").Value = -214298270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6822.425
$ws.Range("I132").Value = 7977
$ws.Range("J132").Value = 3358.7
$ws.Range("K132").Value = 23931
$ws.Range("L132").Value = 10076.1
$ws.Range("M132").Value = -21401
$ws.Range("N132").Value = -15136.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1716.4667
$ws.Range("I122").Value = 1653.6046
$ws.Range("J122").Value = 1875.4706
$ws.Range("K122").Value = 4960.8138
$ws.Range("L122").Value = 5626.4118
$ws.Range("M122").Value = -2510.8138
$ws.Range("N122").Value = -10526.4118
